$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new shared string / cell value in B6
$ws.Range("B6").Value = "sure bitti"

# Update the selection to match the new active cell
$ws.Range("I8").Select() | Out-Null
